$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"
$ws.Range("B3").Value = 1859.8
$ws.Range("D3").Value = "Novo  |  +10mil vendidos"
$ws.Range("E3").Value = "https://www.mercadolivre.com.br/smart-tv-lg-50-4k-uhd-hdr-thinq-ai-pro-wi-fi-bluetooth-alexa-apple-airplay-50tu801c0sa/p/MLB52058084#polycard_client=search-nordic&search_layout=grid&position=1&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4186021237&sid=search"
$ws.Range("F3").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 4
$ws.Range("A4").Value = "Smart TV 4K 50`" LG UHD 50UA75 Processador α7 AI Ger8 4K Super Upscaling Google Cast Alexa Integrado Controle Remoto Padrão WebOS 25"
$ws.Range("B4").Value = 1834.16
$ws.Range("D4").Value = "Novo  |  +1000 vendidos"
$ws.Range("E4").Value = "https://www.mercadolivre.com.br/smart-tv-4k-50-lg-uhd-50ua75-processador-7-ai-ger8-4k-super-upscaling-google-cast-alexa-integrado-controle-remoto-padro-webos-25/p/MLB53201474#polycard_client=search-nordic&search_layout=grid&position=2&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5547463722&sid=search"
$ws.Range("F4").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 5
$ws.Range("A5").Value = "Smart Tv LG Led 55 55ut801c0sa 4k Ultra Hd Thinq Ai Webos 24"
$ws.Range("B5").Value = 2079.44
$ws.Range("D5").Value = "Novo  |  +10mil vendidos"
$ws.Range("E5").Value = "https://www.mercadolivre.com.br/smart-tv-lg-led-55-55ut801c0sa-4k-ultra-hd-thinq-ai-webos-24/p/MLB44672028#polycard_client=search-nordic&search_layout=grid&position=3&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB3966901001&sid=search"
$ws.Range("F5").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 6
$ws.Range("A6").Value = "Smart TV 4K 50`" LG UHD 50UA85 Processador α7 AI Ger8 4K Super Upscaling Google Cast Alexa Integrado Controle AI Smart Magic WebOS 25"
$ws.Range("B6").Value = 1877.8
$ws.Range("D6").Value = "Novo  |  +1000 vendidos"
$ws.Range("E6").Value = "https://www.mercadolivre.com.br/smart-tv-4k-50-lg-uhd-50ua85-processador-7-ai-ger8-4k-super-upscaling-google-cast-alexa-integrado-controle-ai-smart-magic-webos-25/p/MLB51381761#polycard_client=search-nordic&search_layout=grid&position=4&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5474190160&sid=search"
$ws.Range("F6").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 7
$ws.Range("A7").Value = "Smart TV LED 43`" Ultra HD 4K LG 43UR781C0SA.BWZ Preto ThinQ AI 3 HDMI 2 USB Wi-Fi Bluetooth HDR10"
$ws.Range("B7").Value = 1699
$ws.Range("D7").Value = "Novo  |  +50mil vendidos"
$ws.Range("E7").Value = "https://www.mercadolivre.com.br/smart-tv-led-43-ultra-hd-4k-lg-43ur781c0sabwz-preto-thinq-ai-3-hdmi-2-usb-wi-fi-bluetooth-hdr10/p/MLB25800955#polycard_client=search-nordic&search_layout=grid&position=5&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5343939420&sid=search"
$ws.Range("F7").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 8
$ws.Range("A8").Value = "Smart TV Profissional 4K 50`" LG UHD 50AU801 Processador α7 AI Ger8 Super Upscaling Google Cast Alexa Integrado Controle AI Smart Magic WebOS 25"
$ws.Range("B8").Value = 1895.73
$ws.Range("D8").Value = "Novo  |  +25 vendidos"
$ws.Range("E8").Value = "https://www.mercadolivre.com.br/smart-tv-profissional-4k-50-lg-uhd-50au801-processador-7-ai-ger8-super-upscaling-google-cast-alexa-integrado-controle-ai-smart-magic-webos-25/p/MLB61517857#polycard_client=search-nordic&search_layout=grid&position=6&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5889143382&sid=search"
$ws.Range("F8").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 9
$ws.Range("A9").Value = "Smart TV 4K 50`" LG QNED 50QNED82 Processador α7 AI Ger8 Ultra Slim Design Dolby Atmos Alexa Integrado WebOS 25"
$ws.Range("B9").Value = 2881.14
$ws.Range("D9").Value = "Novo  |  +100 vendidos"
$ws.Range("E9").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=7Ioj7EtgdaXlnCSPxzB7%2FrpQaIKTcqWuCIq73sfv25vwB69u1TeHARlWa6lRZhhbGLF9YrZkF1cHNcfhAA%2FgW%2Fjt0ofgR1pWC3ACTp6rzbcqhg2iKKxeXq9JOJyzVCUB6B9DHlPY%2FS4K1l9qfGTP8KgRQ16O%2FD%2FOohaqk%2BkCkZLz79PjkBOhUUrZDup%2FGfUT1sQMbOmapP1Q7L8UAgW6V3b9qJa%2F%2BC7p9XfBr3zoz1KptULIbiQH1HZqZhdQz7M4KLCnusz9yZdlpx6KR%2FOKc0A5uRVPD3JEOzsCcSX04tsoExb4mOKDYRrwthLqqhW3Rc8h74vb6A3MRhzzHAw2qs2ngamWYDNkOS5sIy2boJHqAjcPbFWzNIBDTJpft50Sr71AUoTtFZWOU%2BE45iz28fFBUqGUaBHM8vz514ZmLFk2cpXnPRYp8LeTrWMR4X6r3pXsC74GhvLwRccS%2BlznJLddEYEmRoGeBKUt6P%2BKlCcJEPOrJt4zmRSGy1BgR6%2F54jD06U7Qe87H4tKb%2FxrsOxkB872AdjvmXOPrhomjM%2BpeyRu%2FTjk2WpTAlb01sjQYJ%2F1JqYYPh27zw7LtYuyyKzAtZESBmwqh1zXeri21IcrXiOk%2FUr27kxg6Ae%2B8DTkJQ26A1DboYnM%2FvAPK8xB61GQPh24C0sdeAE9KtBeXECRTXf1kCEXTTV75FrO2GZBTHc%2FAIbQPmkli05ToYc7v8UPDYw%2B9n7atMKhNWLKuu18W0rTAy7S%2Bnn5s%2FK677LsBpaxDt4OvX%2Bn0FQRxeZozswkueJKR6vAm7JagiiuqE3%2FVjr8t%2BYlis0DZtbLO0ap%2B1AHAsDu6vsv74kiszJXNHoiNqHxzm0NYFg6K9JS5HCRv8rh5yMsqGHbBZfkGua%2Fe1OFX7PK4Iu6Pb5lZg95o8TGLvyvFbyv8hDiWAqF967e9ZYPkl0Y0W3arIYP8L36ngDXBNXmtaj%2BUULWX8LYzQ5fnhuUWEVB0iFD0eTDjmatr6TX%2BLH5Me9gY1sbvn2WnN0sEuQ5Tab9CKXK8hqX1C6xkDZEZIQp9L0X%2BLw%3D%3D#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB53712521&backend_model=sup_pads_rars_pdp_pads_up_v2_with_defaults&search_layout=grid&position=7&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4157368403&sid=search"
$ws.Range("F9").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 10
$ws.Range("A10").Value = "Samsung Vision AI TV 65`" NEO QLED Ultra 4K QN70F 2025, Mini LED, Processador com AI, Controle por Gestos, Modo AI, Painel até 144hz"
$ws.Range("B10").Value = 5735.79
$ws.Range("D10").Value = "Novo  |  +50 vendidos"
$ws.Range("E10").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=aBOlOYUIFf1nc2bgqHy%2F%2BI6X1467lyb8z8F3IFzkPS75DBhNkhIclJGOS5MJo45MWECVnvd2c57SdgFzsy%2BqfYw2u9u5Wha9RFIAFT24mz2yLHCn%2FM4MOxj1iUoxH7PWKX2WSRIOzARHEROiMmJpb6ShDX0Zm29EnFLCF90Pvz%2BinkL3bo%2BKrepOvIR7zotlNvZVOr83ebmzr%2Bv9tpumOiVu5gef1YlQZwmFF132dWHYhR1pID%2FxjFbacqJ%2BxgqZDlZcGr0mUUH6ObmkIkcXl1BIfoyU6Mv9vfA3iZ01ziNv2%2FMYHkONVVI2Zea6clDYZpGqyV8vuwIUjyAKYc7uAFMILtFb1aqqKDY0IDCuE%2BzRRtNJR%2BjZtzaGr9V0lknnPx%2Basvj0L1zAV2%2FRhZeiJ%2BboaehU9%2Fr0%2F7%2BOGMOX0qFsjJK%2BSrFasKENCfto8J%2BI%2FIKUjxAjTPmCn%2B%2BtSdCLpHPW8g6ds8HNfxjW4QhzubNlOVhNSq05X%2Fr3OmqTVSPJ2DNtiwLTQFspc2qyAtRcmF1WpmyuGpqCkSbrt3uC1mMdIo%2BKTA7uNVm0kvzmnb0NXIrqIbpm73dQ5qCCXXM1NiaL3W%2BLDEhW2JDimNcccVFTPSD0HuQvAEvTDxPrb%2Fe7vvHDrATem9I4kT%2B4AZmVRMB%2FWcK8ZfGAL5ew0UdqjqbZifLBmyOmCK37OtfVvSn9RZ8DCoHxUzpHERggr3Sb%2Bcnpq6SnYXQCGdLjnOuepYJ7cS3Bbka%2BBUb3LwhA%2BaEBoPbghWEvL6AC0MlTOf%2BeSuaFbMCpIQ2fh1sRVe70EGEP49gqhPrbjaABGLLnJAr6oS4ZhMeo2mV09vcFii5beoPnvqeUh5xJK6gp4oDUxnCIF%2BknR3W2e7WxZESuEMNRbj17D8F3f%2FBuk5iSc5TIrwm%2FG%2BQWeV1B%2FTJVHd%2F4u23VEqwLy8Bb%2B4Gl3V9Dta1Gq9mCkdmKQC9vpB0erIzgqxPJHBj54ddYTeYFceHx%2B%2FtAwYp8L3M47zqtu1LW2NrdoxPbxLroR%2FdiVRRWosh1twViIe%2FE7BEqdBDS48zZKvHTY8jxdP83Ei96Bhr6#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB54054093&backend_model=sup_pads_rars_pdp_pads_up_v2_with_defaults&search_layout=grid&position=8&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5886833650&sid=search"
$ws.Range("F10").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 11
$ws.Range("A11").Value = "Smart TV TCL 43 Polegadas Full HD QLED S5K WiFi Bluetooth Google TV 2 HDMI HDR10 Dolby Audio 43S5K"
$ws.Range("B11").Value = 1376.55
$ws.Range("D11").Value = "Novo  |  +500 vendidos"
$ws.Range("E11").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=pwbkpvO%2Bk45XoXOIIt%2BdwEUZQZnb2xpt48djT5zwsekAiW2PEnmtAEbOxOddhPhfyQchPLOfzQSlENp8NtU14XIUcyo10aaZHXu5iA0XV3ywAqnih9KiHyfKcZvMCfN2KwRiCOCrH6b8YxpXIkgTA5xHaQhVHH%2BKBnAoqJDAsimyVJFKrcuFxXrF1fvp33C%2FvyRmsxgWrUx1tfAb84Jbe1AsZFfVNNqovd5%2B3AvOf2QyJSB8rmSolQ4QtmswFsuNNSfCmVVkp48BUCCcfmIea1rjVW6DTi9Th3L17ZG252%2F2mJIF1NRMrcpj2h%2FsvKNQgR6%2FKghMuehrBcpQYT%2F8TOiv5y14Wy%2BqlEh%2FWlgKJ8is%2BPzetujCfd3PZehgiz%2B1k5nYue9JFNSz4gET%2FsdfhgQh%2BcLD%2F1I7GNKCahGTHK1A2uHczr2uAMZk%2FnJC122W%2FutTttTuKaBUdowqKMnkpD17Ip8bf%2By8FVXPS50PoViS%2BQA5lgURVcwaNA4CTsU%2Bbn7CXPVNUovU7LhCDGnsHcTaTmAidfQRyBcFBD9AifkQkiIDZ5IG1VaZGXOS8c3LBuFZBvhjZ7t0FCyCNhYbt8F2nKtZExg7E8E1PHgrT5CtTdHY3WNWrLeF15mAXElx5u%2Bbf6GWVKHjjwLniDLFYir1WfxoxsLD1LR5a1YkGY38KnKolbaii9Gsmk9w8eIozDTX2e6iWHdFu3Sghs9APxGVjr4ub6p%2FSW93rQApR15UG4%2Bw4ROufCvjPVymfCAC%2B2hGGtizy%2Fb2ZUJcZX5eb%2BVuHXy%2FgB1sJwZp3WBi74D3RyCN7Ns2JIZgp72XybOvNZCkR%2Bh2RkAuQ5SSUKTSNNQ%2ByDoAobSpCqHRJvBTZRHtxgn23GHJBFVhx68fA3e%2Bheu1yn3pcQ9z4TFCLaLxI%2FpTk%2Bp7zLOv5wi57ih739diAs4jOuFho9DphDD4D7KmKYBm7rgmnwC%2FtvY8saZ6mAHRqzY56XWnNVYU4Yi8iyQnx6vxdZPB9BLVG2wUtpQ3qjOidlZKgHqPzkjqXe5a%2Fw%3D%3D#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB49197748&backend_model=sup_pads_rars_pdp_pads_up_v2_with_defaults&search_layout=grid&position=9&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5690651572&sid=search"
$ws.Range("F11").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 12
$ws.Range("A12").Value = "Smart TV de 50`" LG UT801C 50UT801C0SA com tela LED 4K 127/220V - preto"
$ws.Range("B12").Value = 1917.52
$ws.Range("D12").Value = "Novo  |  +10mil vendidos"
$ws.Range("E12").Value = "https://www.mercadolivre.com.br/smart-tv-de-50-lg-ut801c-50ut801c0sa-com-tela-led-4k-127220v-preto/p/MLB45672004#polycard_client=search-nordic&search_layout=grid&position=10&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB3968846129&sid=search"
$ws.Range("F12").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 13
$ws.Range("A13").Value = "Smart TV LG 50`" 4K UHD LED Alexa, ThinQ AI, HDR10, WebOS 24 Controle Smart Magic- 50UT8050"
$ws.Range("B13").Value = 2224.04
$ws.Range("D13").Value = "Novo  |  +10mil vendidos"
$ws.Range("E13").Value = "https://www.mercadolivre.com.br/smart-tv-lg-50-4k-uhd-led-alexa-thinq-ai-hdr10-webos-24-controle-smart-magic-50ut8050/p/MLB39457123#polycard_client=search-nordic&search_layout=grid&position=11&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB3845180283&sid=search"
$ws.Range("F13").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 14
$ws.Range("A14").Value = "Smart Tv 4k 50 LG Uhd 50ut8000 Processador 5 Ger7 Ai Alexa/chromecast Integrado Otimizador De Jogos Webos 24 Com Controle Remoto Padrão"
$ws.Range("B14").Value = 1820.46
$ws.Range("D14").Value = "Novo  |  +1000 vendidos"
$ws.Range("E14").Value = "https://www.mercadolivre.com.br/smart-tv-4k-50-lg-uhd-50ut8000-processador-5-ger7-ai-alexachromecast-integrado-otimizador-de-jogos-webos-24-com-controle-remoto-padro/p/MLB47854356#polycard_client=search-nordic&search_layout=grid&position=12&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4029394209&sid=search"
$ws.Range("F14").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 15
$ws.Range("A15").Value = "Smart TV 4K 50`" LG NanoCell 50NANOA Processador α7 AI Ger8 4K Super Upscaling Google Cast Alexa Integrado Controle AI Smart Magic WebOS 25"
$ws.Range("B15").Value = 2109.6
$ws.Range("D15").Value = "Novo  |  +100 vendidos"
$ws.Range("E15").Value = "https://www.mercadolivre.com.br/smart-tv-4k-50-lg-nanocell-50nanoa-processador-7-ai-ger8-4k-super-upscaling-google-cast-alexa-integrado-controle-ai-smart-magic-webos-25/p/MLB51967528#polycard_client=search-nordic&search_layout=grid&position=13&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4238205859&sid=search"
$ws.Range("F15").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 16
$ws.Range("A16").Value = "Smart TV LG 55UT8050PSA 55`" 4K Ultra HD"
$ws.Range("B16").Value = 2184.05
$ws.Range("D16").Value = "Novo  |  +10mil vendidos"
$ws.Range("E16").Value = "https://www.mercadolivre.com.br/smart-tv-lg-55ut8050psa-55-4k-ultra-hd/p/MLB39457012#polycard_client=search-nordic&search_layout=grid&position=14&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5099316684&sid=search"
$ws.Range("F16").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 17
$ws.Range("A17").Value = "Smart Tv LG 55 55ua8550psa+ubook Kit 10 Livros Digitais Geek"
$ws.Range("B17").Value = 2499
$ws.Range("D17").Value = "Novo  |  2 vendidos"
$ws.Range("E17").Value = "https://www.mercadolivre.com.br/smart-tv-lg-55-55ua8550psaubook-kit-10-livros-digitais-geek/p/MLB2041470807#polycard_client=search-nordic&search_layout=grid&position=15&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4228845243&sid=search"
$ws.Range("F17").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 18
$ws.Range("A18").Value = "Smart TV 4K 65`" LG UHD 65UT8050 Processador α5 Ger7 AI Alexa/Chromecast integrado Otimizador de Jogos webOS 24 Controle Smart Magic"
$ws.Range("B18").Value = 3149.1
$ws.Range("D18").Value = "Novo  |  +5mil vendidos"
$ws.Range("E18").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=HqTEam%2FMmxKZnGSNRAoBnYUMzrwe0wGKrbNYOkhBKG19quOgvzbIQyBsd3%2FFSu4Bc%2BZ3PwXSIiEKH%2ByawXo2gpZ7bAtVRGAMCWjhvlaz1R5xCnwaEDS5YPX1OUlC2iwa7s7g5uvE8NYU2WrWDBQb5i8jfL46Th6zUCAQ7mpqNy3jzzs2mF4IUpyXYZpVHabrKhDtuGNcgYQh8Pp7d03pOypt6RQo1z6fmjO9MQV0sZFlcnAT%2BeGC7Dlr11uJ%2Bo3LqkSWRVEtY4LTfxxGcIW9KzncduTmp%2F0P1LJ%2BR9BhVPxWHaOYencva1Mnp2mOwIacrEF6b8vGYEwuKW0jfALcTCvpmlwcULxyfXweUa7cuBZc75iKvnofwBMunNVfledsxLZbPuDfmkgVibT1H7fyyy9q1nTfxQeCUR2GJY%2FmpxBYWupr4XYrPe%2ByXC4bnQ30nanCC3kU25hE53ZcNUodIUO7EqVhjv2NfIJldgFKIAXhvlOWkGW7an3VLlGULyqQswtTJFChujm0dQSxfr6KlzcklUn5jKyg%2BsM1kvrPC5XiS%2Fif7DuAPwwetVZMd87y6N5%2FnXOQ9tdB2JjzQ0jzHqKyrTQcaVBUV712elt1TtAAQBIm9ofgn50VpugH1pnI8Nfg%2BL%2F4wYF1sH%2F6OqwrpmX%2FXT3%2Bf8DcTaVEsxFEv1o7q3m5NL%2Bfzf8ytcj5TX%2FSNz2u3y7gSQr5JC%2Buzv%2FC9Uu0ty%2BMvnz91nbT5HzwIoVTZeGiUvzZlMrSicCocHo9GsoMZwsmqmeXwc4ulgsOsMZ7lCChkMHbcP%2B2FwjNJ7N9UWbnYdool0I9U390lo3MSkJe4iouot1cZr8ro%2BE4v73%2B4qc8Cta4rlv%2FMw9uo4CHEtcaPa26RCDZltX5k2IJRyp9oHBrZHmAfmkAJm6q3MMmBJvdC%2BH6ghqXz0CyEDLd8tLTsKREBQrMnSpdWkiR4xGnzMeFwBE5zifcS4PBQtGgL%2Bbi4f51mppreusWHMWw8lpuYtGjctGSrWJhEPXM6JXCW%2BMyEcC9oKC0NneXpJdQQh3w#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB39457851&backend_model=search-backend&search_layout=grid&position=16&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB3866900173&sid=search"
$ws.Range("F18").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 19
$ws.Range("A19").Value = "Smart Tv 43'' Semp 43s62 Led Hva Google Tv Hdr Uhd 4k Hdmi"
$ws.Range("B19").Value = 1662
$ws.Range("D19").Value = "Novo  |  +100 vendidos"
$ws.Range("E19").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=WJ5ICkiqgqdIA4SKwM9SOp3ubuCksnye0fvRuZcqBgOwlH6hfdq9F7JJQVg6fxT9OF1gisT%2BxzYWuEl%2BIWNLdL0PgmMXKCN0FfORbY2bV6JDdJ11lSduq6kP5kSlcQVnTsFapJkvSiHphqcrrTYPqxQEL1omszlpAkx%2FREUBj05Igeq8ZoQMvteE2mL%2BavMwnerlX9z5xfPoe4hWZXbEM0bry5bpkCJvfv9aFS%2BgaregRQDUyFp7Rrr5MacqPoFkUJjiHSeOJRV9W5%2BvZFe1iq9zMngHUqlhv4KK%2FVavC4%2FORKaJMNYLtjtGw9I8VF%2B9shde8Pb9gHGp9Qnu1%2BBSQnPhdpZ1Xq%2BAxT%2Frg6tihC05Iw%2B95UjYJ39XC3vo7TisWA1f7ch6kKAKgpiRYWTrUWv2OADC3Y937cujLIu%2FTQt8%2FYkozAKu%2BE6bXMNO7aUegKeOYd4Xw%2F2sVkOEIBMX9FSQ6vzLGPaVj86q2oMm5Med%2FdLLRpmahdIoDKAKw0ssbubrLHHl8OzLip2pZj99wZcqIZXKC1v5mUDxfpkXtQonXxVlZWx%2BXpwrQzK8wIDLs2brYfvns1BKUtNe3kcNUSrnem%2Bt0RInMh2PE7EUkV2qHWWoujnK%2B2XV9nhgrpthqnDdloF8F1BukPFnpwHOSt7fkvKLH5X1edvNFtNK03sPtJKcu8Qz4SR%2Fz0pU90zEucHJ%2FjBAp8XkwchmwHaCA%2F%2FPpAr8EjwpTKjqy39%2FjfUR%2BZAJPfBukwl3feI3bQI3uONe9OksQ2C9e9i7iXhFAKpPgJxWiQi%2FZG7nzazGjW90dLJ27IttB7E1l3u%2F5SgDdHQo7uzhcrUJxaze8z2GalvKefKxO360w4SBnBTKC%2FiyAgv%2BoEF35PpLFbrtOETPCJ0L3ZaaXmrScL48%2Bznvf5wExEVwuiaxK522cssscQ1BHfuOJ0tuoXHywu6MjUD5qzfyaskRW78qEpeMNQLvtJEJi%2ByDPh7kYkeg#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB54075534&backend_model=sup_pads_rars_pdp_pads_up_v2_with_defaults&search_layout=grid&position=17&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5626663944&sid=search"
$ws.Range("F19").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 20
$ws.Range("A20").Value = "Smart Tv 32 Philco Roku Tv Led Dolby Áudio"
$ws.Range("B20").Value = 929.07
$ws.Range("D20").Value = "Novo  |  +10mil vendidos"
$ws.Range("E20").Value = "https://click1.mercadolivre.com.br/mclics/clicks/external/MLB/count?a=NjRH0h6wqKQblNMKYnY80CwlVXWkUjaNud%2BAQ3dPOkxXrvAjNZb1S7l7ujRVVJMLFTUgV7jN157Tl%2FiKbWLNVEdEx8MoCEwOxKvv7UZSNnxIwKsk%2Fbv6MekGEVB%2BhF7jc2O0cxcPTpHMMDXFD6%2BpRh98GeZw%2F5SyqGK%2FzxaLewfQ9EJ7IyIxpP0MDzeRe8EFcG%2FHlNK%2FQecpBT8qieiQoaOBhiXCVxlplnKboFZqVv%2BLPTHNLNcrT7n6ArEn%2BPkwtW2ZiDkDG%2F2KLTlQ7tyHHZFk%2BTLRsuIN%2F99SHJRKPPckXuXeG5NBR%2BdJQnZ9mfR1b7o66Tt%2B%2FyLljbRwo6GHSWmK2om6GXP6dyXbbxBqP9I3T9%2BmWgRfD6FULRQsirVXPi4e4nNFeB%2FSgU8Fxj5ZgKw6Y1LyzYXHNRh2qyD3FFv4Vy7p40tTuY5HRg9G3pjaJOoddGUv2cE1UC6iKPhOhAipI0mJis8Rxu%2FwPT7xQ7wNn2To2SeFqFZWHPnxQFdPCsxMVQgignIV5jNXVSakd5ryTVot794dip9Tshv1P%2BuZJ64lXVb5nBuCRkKq4IYL2lkz5CG%2BZh50qqqGGne5edoAzYjIe9a5IENftGRJZo0knXg2q%2BTi7euyX7FZleePg3P%2BIjSCrC6Ztj93CepdfUixGjMUf5pb7v24pXDZyUzBtnky9nVJyTUO0rqBRASinehDjKPH1KGemexu1zWptPIo%2FX9lZybFg7ixBbXDKi%2FWxM7bdeQKqQVZZNmUK2u7NREaEPqOJDnSFHfq8GAt%2FHzkFnzVEVStmPgvxR9eQrHVHEr3pO%2BJoSc7%2FCoMK7tZ3XU0zRrmxh6A58HLKJsFKAe0leUjqr%2FXe25x6a7LLLvV8VrrLLm%2FZvkoJV9CX4kn5gVFPZOF466cYDQ9AI6PBRC2X%2FtXoNxn8DeL7SkZFoY2UrO%2Bpxx7Ip%2BSLvr9CUwTV%2FRLcUtXUIqAUiI%3D#polycard_client=search-nordic&is_advertising=true&searchVariation=MLB23723335&backend_model=sup_pads_rars_pdp_pads_up_v2_with_defaults&search_layout=grid&position=18&type=pad&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB4020502875&sid=search"
$ws.Range("F20").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 21
$ws.Range("A21").Value = "Tv 50  Lg 4k Uhd Smart Inteligencia Artificial Thinq - 50uq8"
$ws.Range("B21").Value = 1963.66
$ws.Range("D21").Value = "Novo  |  4 vendidos"
$ws.Range("E21").Value = "https://www.mercadolivre.com.br/tv-50--lg-4k-uhd-smart-inteligencia-artificial-thinq--50uq8/up/MLBU3139674981#polycard_client=search-nordic&search_layout=grid&position=19&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5363831990&sid=search"
$ws.Range("F21").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"

# Row 22
$ws.Range("A22").Value = "Smart Tv Lg 50  Smart 4k Uhd Nanocel 50nano80tsa.bwz"
$ws.Range("B22").Value = 3223.75
$ws.Range("D22").Value = "Novo  |  +100 vendidos"
$ws.Range("E22").Value = "https://www.mercadolivre.com.br/smart-tv-lg-50--smart-4k-uhd-nanocel-50nano80tsabwz/up/MLBU3068438982#polycard_client=search-nordic&search_layout=grid&position=20&type=product&tracking_id=3d2649de-c3ae-4c5c-8643-9fefddd0e4d2&wid=MLB5316304420&sid=search"
$ws.Range("F22").Value = "Smart Tv LG 50 4k Uhd Hdr Thinq Ai Pro Wi-fi Bluetooth Alexa Apple Airplay - 50tu801c0sa"
